$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

# Corrected Observed data (pivot table in source spreadsheet was not updated)
# The following cell values in column E are updated to match the corrected pivot output.

$ws.Range("E11").Value = 2.027840351393797
$ws.Range("E12").Value = 3.353731066541702
$ws.Range("E13").Value = 4.719348659003832
$ws.Range("E14").Value = 5.777999158767037
$ws.Range("E15").Value = 6.591597796143251
$ws.Range("E16").Value = 7.285626312069247
$ws.Range("E17").Value = 8.164738253139552
$ws.Range("E20").Value = 2.250978473581213
$ws.Range("E21").Value = 3.522537971582558
$ws.Range("E22").Value = 4.602474226804123
$ws.Range("E23").Value = 5.662004914715236
$ws.Range("E24").Value = 6.393929684556216
$ws.Range("E25").Value = 7.168551106121715
$ws.Range("E26").Value = 8.016140222803642
$ws.Range("E28").Value = 3.120248166580074
$ws.Range("E29").Value = 4.280431539301513
$ws.Range("E30").Value = 5.35892959870313
$ws.Range("E31").Value = 5.970673860292429
$ws.Range("E32").Value = 7.333638743631557
$ws.Range("E33").Value = 7.938889123888614
$ws.Range("E34").Value = 8.925237610053305
$ws.Range("E35").Value = 9.470722853535353
$ws.Range("E42").Value = 2.18912447257384
$ws.Range("E43").Value = 3.429245283018868
$ws.Range("E44").Value = 4.936791612694301
$ws.Range("E45").Value = 6.249056186868686
$ws.Range("E46").Value = 7.176190476190476
$ws.Range("E47").Value = 8.423809523809524
$ws.Range("E50").Value = 2.052940447297685
$ws.Range("E51").Value = 3.501275510204081
$ws.Range("E52").Value = 4.59433962264151
$ws.Range("E53").Value = 5.534394589244473
$ws.Range("E54").Value = 6.739725753999156
$ws.Range("E55").Value = 6.990567609181475
$ws.Range("E56").Value = 7.734101020675215
$ws.Range("E73").Value = 2.781931878658861
$ws.Range("E74").Value = 3.690562652035528
$ws.Range("E75").Value = 4.36478244834977
$ws.Range("E76").Value = 5.082653952716075
$ws.Range("E77").Value = 6.238632385698204
$ws.Range("E78").Value = 6.801271651422589
$ws.Range("E79").Value = 7.622494103165831
$ws.Range("E80").Value = 8.399167381292468
$ws.Range("E81").Value = 8.986087073665527
$ws.Range("E82").Value = 10.08956131078224
$ws.Range("E83").Value = 10.53762678632404
$ws.Range("E84").Value = 11.88651674985288
$ws.Range("E87").Value = 2.860614359812277
$ws.Range("E88").Value = 3.926960478652095
$ws.Range("E89").Value = 4.303664908238337
$ws.Range("E90").Value = 5.3374275877576
$ws.Range("E91").Value = 6.229974485520966
$ws.Range("E92").Value = 6.918292467890037
$ws.Range("E93").Value = 7.557093822446562
$ws.Range("E94").Value = 8.434939644520483
$ws.Range("E95").Value = 8.785902197363226
$ws.Range("E96").Value = 10.14013495013843
$ws.Range("E97").Value = 10.34354844003215
$ws.Range("E101").Value = 3.064375519026693
$ws.Range("E102").Value = 4.126303013883644
$ws.Range("E103").Value = 4.859282316201672
$ws.Range("E104").Value = 5.636901512027431
$ws.Range("E105").Value = 6.894152987959807
$ws.Range("E106").Value = 7.539066891512086
$ws.Range("E107").Value = 8.636567773124037
$ws.Range("E108").Value = 9.180154667853566
$ws.Range("E109").Value = 10.75585572116479
$ws.Range("E110").Value = 11.48925659221426
$ws.Range("E111").Value = 12.21963457047641
$ws.Range("E112").Value = 13.60035169019963
$ws.Range("E113").Value = 13.41413611575108
$ws.Range("E114").Value = 15.89130434782609
$ws.Range("E119").Value = 3.58272414477179
$ws.Range("E120").Value = 3.63905043229178
$ws.Range("E121").Value = 4.874789756085082
$ws.Range("E122").Value = 6.830815018315017
$ws.Range("E123").Value = 8.518928004677191
$ws.Range("E132").Value = 2.85531113000544
$ws.Range("E133").Value = 4.008620918892745
$ws.Range("E134").Value = 4.683063969502782
$ws.Range("E135").Value = 5.289699321621002
$ws.Range("E136").Value = 6.461893579972173
$ws.Range("E137").Value = 6.940492443564277
$ws.Range("E138").Value = 7.644492906854274
$ws.Range("E139").Value = 8.486113788677581
$ws.Range("E140").Value = 9.0155232076504
$ws.Range("E141").Value = 10.37313875182323
$ws.Range("E142").Value = 10.19104676783664
$ws.Range("E154").Value = 4.349443790849673
$ws.Range("E157").Value = 3.140414263801046
$ws.Range("E158").Value = 4.728264506576705
$ws.Range("E159").Value = 5.045190028431217
$ws.Range("E160").Value = 6.439297583915042
$ws.Range("E161").Value = 7.619329006326453
$ws.Range("E163").Value = 3.323946706887883
$ws.Range("E164").Value = 4.402277807138626
$ws.Range("E165").Value = 5.256666195998101
$ws.Range("E166").Value = 6.052695606280512
$ws.Range("E167").Value = 6.565958979494665
$ws.Range("E169").Value = 2.915552044666688
$ws.Range("E170").Value = 4.743010508567753
$ws.Range("E171").Value = 5.832747363134104
$ws.Range("E173").Value = 2.815234522477118
$ws.Range("E174").Value = 4.322749860746435
$ws.Range("E175").Value = 5.056834744225061
$ws.Range("E176").Value = 5.880703114774386
$ws.Range("E177").Value = 6.746069360415887
$ws.Range("E192").Value = 1.181641482857677
$ws.Range("E193").Value = 2.707680336754799
$ws.Range("E194").Value = 3.928445747800587
$ws.Range("E195").Value = 5.312150736911827
$ws.Range("E196").Value = 7.122617526006056
$ws.Range("E198").Value = 1.472210591315308
$ws.Range("E199").Value = 2.737906131181342
$ws.Range("E200").Value = 4.815872964956653
$ws.Range("E201").Value = 4.83037663219835
$ws.Range("E202").Value = 5.775715284969713
$ws.Range("E203").Value = 6.674429117465286
$ws.Range("E204").Value = 7.874732263606125
$ws.Range("E205").Value = 8.906989015174556
$ws.Range("E206").Value = 13.02280991735537
$ws.Range("E207").Value = 13.25
